$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..25

$Bvals = @(18.70775606929134,18.06538492518789,17.65811965355348,17.4891259129981,17.46088790044268,17.65585253252852,18.4890469186827,20.01320908473959,21.05701976514742,21.51380596167153,21.68407379001475,21.64752559947074,21.52786884448189,21.45421999570624,21.02679464055362,20.7598761672053,20.60465862850933,20.55181749225098,20.78846617878585,21.56308922698408,22.05352849658974,21.79325283714827,20.7755461158609,19.61363986954485)
$Barr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Barr[$i,0] = $Bvals[$i] }
$ws.Range("B2:B25").Value = $Barr

$Dvals = @(10.82657036259021,10.87170441272268,10.90111379709423,10.91352611259641,10.91561303443253,10.90127946056294,10.84178102669226,10.73852128348972,10.67076954090424,10.64169533674341,10.63093577382412,10.63324192476209,10.64080513174978,10.64547037156794,10.67270467141645,10.68985866729164,10.69988962623453,10.70331420664438,10.68801558354703,10.63857685127335,10.60772376055966,10.62405753267493,10.68884831553191,10.76502630264996)
$Darr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Darr[$i,0] = $Dvals[$i] }
$ws.Range("D2:D25").Value = $Darr

$Evals = @(16.32468574740047,16.31585396372092,16.31316344445973,16.31275845268641,16.31273308608356,16.31315517753044,16.32107618238133,16.35808499281442,16.39807871914935,16.41898550219924,16.42728634936991,16.42548164236847,16.41966076093768,16.41614510695641,16.39676648209834,16.38556891757504,16.37938417055913,16.37733424156881,16.38673448813831,16.42136012742528,16.44622454577399,16.43275162489206,16.38620674548791,16.34580262907238)
$Earr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Earr[$i,0] = $Evals[$i] }
$ws.Range("E2:E25").Value = $Earr

$Fvals = @(26.70272535668428,26.85249988212857,26.9530544213713,26.99618217779052,27.00347313842921,26.95362736102987,26.75257848181421,26.4269270454892,26.23007416933886,26.14985052596311,26.12082292268147,26.12701427211494,26.14743526595296,26.16012000357787,26.23550551629593,26.28414801922627,26.31300307262953,26.32292323611027,26.27887908242513,26.1414003683343,26.059431452426,26.10245529109296,26.28125839775,26.5076195340033)
$Farr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Farr[$i,0] = $Fvals[$i] }
$ws.Range("F2:F25").Value = $Farr

$Gvals = @(24.63915112551384,24.63960521028037,24.65149887786935,24.65925595097826,24.66071941176595,24.65159172503731,24.63689188407213,24.70051392934641,24.80375789692529,24.86294463460301,24.88710410850243,24.88182345587721,24.86489731595789,24.85475666450632,24.80013499879965,24.76975006884538,24.75342459578173,24.7480950378889,24.77286552252796,24.86982163076693,24.94336179127986,24.90318560663339,24.77145346408763,24.67337556756481)
$Garr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Garr[$i,0] = $Gvals[$i] }
$ws.Range("G2:G25").Value = $Garr

$Hvals = @(13.16128777008751,13.21168097263219,13.24541814736251,13.25986805021862,13.26230979196842,13.24561018418082,13.17808239744336,13.06789626512253,13.00057092771923,12.97291737817377,12.9628744662788,12.96501829085106,12.97208253855179,12.97646548784768,13.00243809538209,13.01913390587056,13.02901665576006,13.03241079234173,13.01732764412157,12.96999594720324,12.94156243730262,12.95650870060356,13.01814337017237,13.09531642488863)
$Harr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Harr[$i,0] = $Hvals[$i] }
$ws.Range("H2:H25").Value = $Harr

$Ivals = @(23.50926219041557,23.64052696443126,23.7281604774499,23.76563676673318,23.77196614399147,23.72865875554927,23.55305919658823,23.26474061686395,23.08734924054001,23.01418732231462,22.98757134151546,22.99325505582786,23.01197575627491,23.02358468741903,23.09228261702975,23.13636014521146,23.16242121283255,23.17136662127282,23.13159461560708,23.00644743840155,22.9310056952551,22.97068762935984,23.13374686821376,23.33671131030421)
$Iarr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Iarr[$i,0] = $Ivals[$i] }
$ws.Range("I2:I25").Value = $Iarr

$Jvals = @(11.28474219165262,11.26940785471164,11.26138734416448,11.25847308199348,11.25801065826575,11.26134660303117,11.27916750882884,11.3250261288226,11.36515443178478,11.38475825296101,11.39237131842514,11.39072335252735,11.38538081235479,11.38213289797382,11.36390004519433,11.35305720226837,11.34694796223969,11.3449014731155,11.35419830089774,11.38694493881814,11.40944926661946,11.39733892534516,11.35368202218784,11.31147366449954)
$Jarr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Jarr[$i,0] = $Jvals[$i] }
$ws.Range("J2:J25").Value = $Jarr

$Lvals = @(12.26724367131262,11.88898716869377,11.64972163025492,11.55057004766459,11.5340099678919,11.6483909545483,12.13834252490647,13.03893908894975,13.65851937273964,13.93028977418169,14.0316845604312,14.00991594634524,13.93866239152187,13.89481782182287,13.6405495999095,13.48193073958367,13.38975222446834,13.35838207730798,13.4989143339931,13.95963301985288,14.25186008363768,14.09672633248067,13.49123911812388,12.80234168955984)
$Larr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Larr[$i,0] = $Lvals[$i] }
$ws.Range("L2:L25").Value = $Larr

$Ovals = @(19.49656810137638,19.56541487738021,19.61350250402279,19.6345547296217,19.63813819479191,19.61378053434332,19.519095834538,19.37982649810452,19.30613707003015,19.27889634161303,19.26948890284211,19.27147450573429,19.27810416147508,19.28228340590319,19.30804412844203,19.32545978466487,19.33606774112684,19.33976075137587,19.32354466264879,19.27613219073824,19.25044004254257,19.26366648253703,19.32440863483404,19.41249773151413)
$Oarr = New-Object "object[,]" 24,1
for ($i = 0; $i -lt 24; $i++) { $Oarr[$i,0] = $Ovals[$i] }
$ws.Range("O2:O25").Value = $Oarr

Write-Host "Applied 380 kV case values"
